$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.29), maa://25390 (95.95), maa://36681 (87.01)'
$ws.Range('L2').Value = '*maa://24633 (56.05), *maa://30515 (69.31), *maa://34787 (72.86), ***maa://20792 (11.93), maa://39402 (88.89), ***maa://29083 (27.78)'
$ws.Range('T2').Value = 'maa://22742 (91.08), *maa://20791 (63.38)'
$ws.Range('AB2').Value = 'maa://21246 (91.32), maa://36684 (95.79), ***maa://22731 (6.67)'
$ws.Range('P3').Value = 'maa://21249 (94.22), maa://26254 (96.15)'
$ws.Range('T3').Value = 'maa://24617 (88.99), **maa://20790 (44.78), ***maa://37170 (18.97)'
$ws.Range('D4').Value = 'maa://24632 (93.42), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range('T4').Value = 'maa://32509 (98.04), maa://27295 (83.87), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range('T6').Value = 'maa://37411 (85.71)'
$ws.Range('AF6').Value = '*maa://33152 (61.36), ***maa://22770 (27.27)'
$ws.Range('X7').Value = 'maa://22399 (95.17), *maa://22758 (73.77)'
$ws.Range('AF7').Value = '*maa://26191 (68.75), *maa://36671 (69.39), *maa://42530 (64.29)'
$ws.Range('A8').Value = '更新日期：2025.01.03 13:17:45'
$ws.Range('X8').Value = 'maa://21411 (95.94)'
$ws.Range('X9').Value = 'maa://26223 (97.5)'
$ws.Range('AB9').Value = 'maa://28711 (87.38), ***maa://22740 (5.77), **maa://39938 (46.15), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (94.12)'
$ws.Range('D10').Value = '***maa://25695 (19.34), **maa://32237 (40.48), ***maa://34206 (21.74), ***maa://39951 (16.28), ***maa://39243 (28.57)'
$ws.Range('D11').Value = 'maa://36707 (99.42)'
$ws.Range('H13').Value = '*maa://21248 (74.89), **maa://22728 (47.73)'
$ws.Range('P13').Value = 'maa://22676 (91.96), *maa://22583 (75.38), *maa://22500 (57.78)'
$ws.Range('L14').Value = 'maa://26245 (96.55), maa://21288 (96.3), maa://39841 (95.18), maa://36682 (97.37)'
$ws.Range('C15').Value = "'5"
$ws.Range('D15').Value = '*maa://22743 (77.44), maa://22734 (84.03), *maa://30808 (65.08), **maa://36048 (32.61), maa://45058 (100.0)'
$ws.Range('P15').Value = 'maa://24762 (90.38), *maa://22727 (70.0)'
$ws.Range('AF15').Value = 'maa://21364 (80.97), *maa://22766 (70.27), *maa://36666 (78.82)'
$ws.Range('D16').Value = 'maa://21441 (96.35), maa://36679 (93.18), maa://37650 (96.97)'
$ws.Range('AB16').Value = 'maa://26228 (95.45)'
$ws.Range('AF16').Value = '*maa://23911 (65.69), maa://27755 (92.77)'
$ws.Range('D17').Value = 'maa://21624 (83.33)'
$ws.Range('L17').Value = '*maa://21679 (76.92)'
$ws.Range('D18').Value = 'maa://24570 (97.13)'
$ws.Range('AF19').Value = '*maa://21663 (61.76)'
$ws.Range('H20').Value = 'maa://22864 (89.4)'
$ws.Range('L20').Value = 'maa://41331 (85.29)'
$ws.Range('L23').Value = 'maa://39756 (94.34), maa://39875 (93.85)'
$ws.Range('P23').Value = 'maa://30587 (92.02), *maa://29748 (75.78), ***maa://29785 (16.42), *maa://37566 (68.97)'
$ws.Range('AF24').Value = 'maa://22523 (85.71), maa://36672 (80.77), maa://29910 (92.73), **maa://21440 (34.55)'
$ws.Range('D25').Value = 'maa://29753 (94.9)'
$ws.Range('G25').Value = "'4"
$ws.Range('H25').Value = '*maa://29063 (73.83), *maa://25311 (73.53), ***maa://22725 (4.84), maa://45047 (100.0)'
$ws.Range('P25').Value = 'maa://24382 (93.33)'
$ws.Range('AB26').Value = 'maa://42235 (93.33)'
$ws.Range('H27').Value = '**maa://21283 (48.0), maa://34494 (96.67), *maa://39601 (76.47), **maa://36665 (44.44)'
$ws.Range('D28').Value = 'maa://24465 (90.78), maa://25725 (83.53)'
$ws.Range('T28').Value = 'maa://23263 (95.0), *maa://29765 (60.53)'
$ws.Range('X28').Value = 'maa://39929 (89.75), ***maa://39723 (14.29), maa://41749 (91.53)'
$ws.Range('AF29').Value = '*maa://24080 (69.05), ***maa://34960 (8.33), *maa://42865 (78.95)'
$ws.Range('X30').Value = 'maa://39477 (86.67)'
$ws.Range('AB30').Value = 'maa://42979 (96.46), maa://45045 (100.0)'
$ws.Range('L31').Value = 'maa://35926 (93.7), maa://36258 (83.67), *maa://43904 (80.0)'
$ws.Range('H32').Value = 'maa://21895 (97.33), maa://36667 (98.41), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('P33').Value = '*maa://21956 (79.71), *maa://22730 (79.31)'
$ws.Range('T34').Value = 'maa://24526 (93.6)'
$ws.Range('L35').Value = 'maa://41296 (96.72)'
$ws.Range('AF38').Value = 'maa://36697 (86.34)'
$ws.Range('H39').Value = 'maa://25199 (84.82), maa://36670 (87.78), maa://30434 (89.39), ***maa://25036 (16.0), *maa://44165 (66.67), maa://45059 (100.0)'
$ws.Range('P39').Value = 'maa://24709 (91.47)'
$ws.Range('P40').Value = 'maa://23278 (96.0), maa://21386 (95.74), maa://36664 (90.74)'
$ws.Range('P41').Value = '**maa://35616 (38.24), maa://43177 (86.67)'
$ws.Range('H43').Value = 'maa://22525 (92.7), maa://21284 (83.72)'
$ws.Range('H44').Value = 'maa://29768 (97.85), maa://27728 (96.0)'
$ws.Range('H45').Value = 'maa://21229 (84.86), maa://30807 (95.52), *maa://22767 (55.0), ***maa://20796 (13.79), *maa://42459 (77.78)'
$ws.Range('H53').Value = 'maa://32534 (93.58), **maa://32434 (34.78)'
$ws.Range('H55').Value = 'maa://32532 (91.95)'
$ws.Range('H59').Value = 'maa://27746 (83.02), maa://31270 (94.87)'
$ws.Range('H62').Value = 'maa://42981 (96.67), maa://43903 (100.0)'
